$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.662.48"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.493.68"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.61%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.01"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.95"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.518"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.492.63"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("D10").ClearFormats()
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.09"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.17"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.957.96"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000175"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.681.06"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.519.68"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.22"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -6.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.56"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "342.40"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.18"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.53"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.49%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "68.43"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.91"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.45%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0974"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.13"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "519.61"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.30"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.32%  "
$ws.Range("E34").Value = "  -4.83%  "
$ws.Range("E35").Value = "  -4.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.40"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.44"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.42"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.29"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.353"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.75"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.03"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.43"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "146.14"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.552"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.32%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0276"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.76%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.68"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.70"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0750"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.51%  "
